$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.511.26"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.952.12"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.611"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.50"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.373"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0790"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.73%  "
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "2.242.19"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.818"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").Value = "1.959.42"
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("D18").Value = "36.318.91"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("D20").Value = "0.0₃0845"
$ws.Range("E20").Value = "  -3.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "226.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.136"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "159.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.119"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0608"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.55%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.22"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.98%  "
$ws.Range("E38").Value = "  -0.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -12.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0969"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.16"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0208"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.50%  "
$ws.Range("D45").Value = "1.352.88"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.59%  "
$ws.Range("B50").Value = "FTXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +13.14%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.132.48"
$ws.Range("E51").Value = "  -1.02%  "
